# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets (and the Overview
# roll-up, which shares the same "Status" string) to reflect that the
# handoff files have now been handed back.

$wb = $excel.ActiveWorkbook

$targetFileName = "1970cc51-f781-4ac3-a795-b141092e76e7.md"
$targetFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2360f17b8afdba90f1612e44f2232de927814787/e2e/1970cc51-f781-4ac3-a795-b141092e76e7.md"

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        XlfFile = "1970cc51-f781-4ac3-a795-b141092e76e7.3c733eade4c23280b5a08d5b6a82b6588d8771c4.zh-cn.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b4d759889dbfb490fc2ac315fa7fe413ab42f70/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/1970cc51-f781-4ac3-a795-b141092e76e7.3c733eade4c23280b5a08d5b6a82b6588d8771c4.zh-cn.xlf"
        HandbackTime = "2016-01-07 08:37:36"
    },
    @{
        Name = "de-de"
        XlfFile = "1970cc51-f781-4ac3-a795-b141092e76e7.3c733eade4c23280b5a08d5b6a82b6588d8771c4.de-de.xlf"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a54a902ec79036b2f93fc876ab19d39ed9fcad7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/1970cc51-f781-4ac3-a795-b141092e76e7.3c733eade4c23280b5a08d5b6a82b6588d8771c4.de-de.xlf"
        HandbackTime = "2016-01-07 08:37:08"
    }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status moves from "Not yet handed off" to "Handed back".
    $ws.Range("B2").Value = "Handed back"

    # Latest Target File (E2) and Latest Handback File (F2) now get filled
    # in with the same file / xlf that were already referenced in A2 / C2,
    # complete with hyperlinks matching those cells.
    $ws.Hyperlinks.Add($ws.Range("E2"), $targetFileUrl, [Type]::Missing, [Type]::Missing, $targetFileName)
    $ws.Range("E2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("F2"), $info.XlfUrl, [Type]::Missing, [Type]::Missing, $info.XlfFile)
    $ws.Range("F2").Style = "HyperLink"

    # Latest Handback DateTime (G2) is now stamped with the handback time.
    $ws.Range("G2").Value = $info.HandbackTime
}

# The Overview sheet's Status columns (B2/C2) mirror the same shared string,
# so update them to match as well.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back"
$overview.Range("C2").Value = "Handed back"
